# Weekly update: add a new "Papa" (potato) price record for
# "Macroferia Regional de Talca" as the first row of that variety's block
# (row 469), shifting the existing rows (469-497) down by one (470-498).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 469; Excel shifts rows 469-497
# down to 470-498 and carries the formatting (incl. the date style on
# column D) down from the row above, just like a manual "Insert Row".
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A469").Value = 5
$ws.Range("B469").Value = "Macroferia Regional de Talca"
$ws.Range("C469").Value = "Maule"
$ws.Range("D469").Value = 44714
$ws.Range("E469").Value = 7
$ws.Range("F469").Value = 100114001
$ws.Range("G469").Value = "Papa"
$ws.Range("H469").Value = "Asterix"
$ws.Range("I469").Value = "1a (cosecha)"
$ws.Range("J469").Value = 1200
$ws.Range("K469").Value = 7000
$ws.Range("L469").Value = 7000
$ws.Range("M469").Value = 7000
$ws.Range("N469").Value = "`$/saco 25 kilos"
$ws.Range("O469").Value = "Región de Los Lagos"
$ws.Range("P469").Value = 280
$ws.Range("Q469").Value = 25
$ws.Range("R469").Value = "Hortaliza"
